$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update praclen (column I) for existing rows 2-4 from 4 to 5
$ws.Range("I2").Value = 5
$ws.Range("I3").Value = 5
$ws.Range("I4").Value = 5

# Replace row 5 data with the new trial, and add a new row 6 with the old row 5's data
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 61
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = "train_dim1_1"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 66
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = "train_dim1_1"

$ws.Range("K17").Select()
